# Update runs/balls/fours/sixes stats for Rohit Sharma (c) innings rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "68"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "4"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "0"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "0"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "9"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "8"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2"

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "8"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "8"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "1"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "4"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "0"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "80"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "6"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "70"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "8"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "3"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "12"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "10"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "0"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "35"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "36"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "35"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "3"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "6"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "1"
